$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Aslib Proceedings', 1),
    @('British Food Journal', 2),
    @('Contributions to economic analysis', 1),
    @('Data Technologies and Applications', 2),
    @('Disaster Prevention and Management An International Journal', 1),
    @('Education + Training', 2),
    @('Employee Relations', 1),
    @('Engineering Computations', 2),
    @('Environmental Management and Health', 1),
    @('Equal Opportunities International', 1),
    @('European Journal of Innovation Management', 1),
    @('European Journal of Management and Business Economics', 1),
    @('European Journal of Marketing', 2),
    @('Evidence-based HRM a Global Forum for Empirical Scholarship', 1),
    @('foresight', 2),
    @('Gender in Management An International Journal', 1),
    @('Industrial and Commercial Training', 1),
    @('Industrial Lubrication and Tribology', 1),
    @('Industrial Management & Data Systems', 2),
    @('Industrial Robot the international journal of robotics research and application', 1),
    @('Info', 1),
    @('Information and Learning Sciences', 1),
    @('Interactive Technology and Smart Education', 1),
    @('International Journal of Architectural Research Archnet-IJAR', 1),
    @('International Journal of Climate Change Strategies and Management', 2),
    @('International Journal of Conflict Management', 1),
    @('International Journal of Contemporary Hospitality Management', 4),
    @('International Journal of Disaster Resilience in the Built Environment', 2),
    @('International Journal of Educational Management', 7),
    @('International Journal of Entrepreneurial Behaviour & Research', 1),
    @('International Journal of Gender and Entrepreneurship', 4),
    @('International Journal of Information and Learning Technology', 1),
    @('International Journal of Manpower', 1),
    @('International Journal of Operations & Production Management', 1),
    @('International Journal of Physical Distribution & Logistics Management', 3),
    @('International Journal of Social Economics', 4),
    @('International Journal of Sociology and Social Policy', 3),
    @('International Journal of Sustainability in Higher Education', 4),
    @('International Journal of Wine Business Research', 2),
    @('International perspectives on inclusive education', 1),
    @('Internet Research', 2),
    @('Journal for Multicultural Education', 1),
    @('Journal of Aggression Conflict and Peace Research', 2),
    @('Journal of Agribusiness in Developing and Emerging Economies', 2),
    @('Journal of Business and Industrial Marketing', 1),
    @('Journal of Communication Management', 1),
    @('Journal of Consumer Marketing', 2),
    @('Journal of Documentation', 1),
    @('Journal of Economic Studies', 1),
    @('Journal of Educational Administration', 10),
    @('Journal of European Industrial Training', 1),
    @('Journal of Facilities Management', 1),
    @('Journal of Fashion Marketing and Management', 1),
    @('Journal of Health Organization and Management', 3),
    @('Journal of International Education in Business', 1),
    @('Journal of Management Development', 2),
    @('Journal of Managerial Psychology', 3),
    @('Journal of Manufacturing Technology Management', 1),
    @('Journal of Organizational Change Management', 1),
    @('Journal of service management', 4),
    @('Journal of Services Marketing', 1),
    @('Journal of Small Business and Enterprise Development', 2),
    @('Journal of Tourism Futures', 1),
    @('Journal of Workplace Learning', 1),
    @('Kybernetes', 6),
    @('Leadership & Organization Development Journal', 4),
    @('Library Hi Tech', 1),
    @('Library Management', 1),
    @('Library Review', 3),
    @('Literacy research, practice and evaluation', 1),
    @('Management Decision', 3),
    @('Management of Environmental Quality An International Journal', 4),
    @('Management Research News', 1),
    @('Managing Service Quality', 1),
    @('Microelectronics International', 1),
    @('New England journal of entrepreneurship', 1),
    @('Nutrition & Food Science', 2),
    @('On the Horizon The International Journal of Learning Futures', 1),
    @('Organization Management Journal', 1),
    @('Personnel Review', 3),
    @('Policing An International Journal', 3),
    @('Qualitative Market Research An International Journal', 2),
    @('Qualitative Research in Accounting & Management', 1),
    @('Qualitative Research in Organizations and Management An International Journal', 1),
    @('Rapid Prototyping Journal', 11),
    @('Reference Reviews', 13),
    @('Reference Services Review', 2),
    @('Research in community and mental health', 1),
    @('Research in experimental economics', 6),
    @('Research in social movements, conflicts and change', 3),
    @('Research in social problems and public policy', 2),
    @('Research in the history of economic thought and methodology', 1),
    @('Research in the sociology of health care', 1),
    @('Research in the sociology of organizations', 9),
    @('Review of Accounting and Finance', 1),
    @('Review of marketing research', 1),
    @('Sensor Review', 9),
    @('Smart and Sustainable Built Environment', 1),
    @('Social Responsibility Journal', 1),
    @('Society and Business Review', 1),
    @('Strategic HR Review', 1),
    @('Strategy and Leadership', 2),
    @('Studies in Economics and Finance', 1),
    @('Studies in media and communications', 2),
    @('Technology, innovation, entrepreneurship and competitive strategy', 2),
    @('The Electronic Library', 1),
    @('Tizard Learning Disability Review', 2),
    @('VINE', 1),
    @('Work Study', 1),
    @('Young Consumers Insight and Ideas for Responsible Marketers', 2)
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove now-obsolete trailing rows (sheet previously ended at row 121, now ends at row 116)
$ws.Range("A117:B121").ClearContents()
